# Global Update et correction addTraceBEGIN et END
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# The reference list in column A grows from 2 entries (+1 trailer) to
# 8 entries (+1 trailer). Write the new values directly; D4 (already
# present, empty, formatted as text) is left untouched in place.
$ws.Range("A2").Value = "AD.SEC.001.FON.02"
$ws.Range("A3").Value = "AD.SEC.001.FON.01"
$ws.Range("A4").Value = "AD.SEC.001.FON.03"
$ws.Range("A5").Value = "AD.DEP.001.FON.01"
$ws.Range("A6").Value = "RO.ACT"
$ws.Range("A7").Value = "RO.FOU"
$ws.Range("A8").Value = "MP.CPT"
$ws.Range("A9").Value = "RT.ART"
$ws.Range("A10").Value = "AD.SEC.014.FON.01"

# A2:A9 all share the existing "text" cell format (same as the original
# A2 / D4 cells); A10 keeps the General format the trailing row had.
$ws.Range("A2:A9").NumberFormat = "@"

# Selection moves from C8 to C7, and the used range grows to A1:D10.
$ws.Range("C7").Select()

$wb.Save()
